$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Set header values for new columns
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Copy formatting from an existing header cell (AC1) to new header cells
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)  # xlPasteFormats

# Fill season record (Wins/Losses/Ties) for each player row
$ws.Range("AD2:AD42").Value = 91
$ws.Range("AE2:AE42").Value = 71
$ws.Range("AF2:AF42").Value = 0

Write-Host "Done"
